$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.450.13'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '3.608.78'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.40'
$ws.Range('E5').Value = '  -1.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '190.06'
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('D7').Value = '3.605.57'
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('E8').Value = '  -1.92%  '
$ws.Range('E10').Value = '  +3.91%  '
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.11'
$ws.Range('E12').Value = '  -3.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000313'
$ws.Range('E13').Value = '  +8.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.70'
$ws.Range('E14').Value = '  -2.16%  '
$ws.Range('D15').Value = '4.189.72'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.86'
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('D17').Value = '3.608.77'
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').Value = '70.404.90'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.64'
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '493.26'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '19.35'
$ws.Range('E23').Value = '  -1.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.92'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '96.99'
$ws.Range('E25').Value = '  +6.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.37'
$ws.Range('E26').Value = '  -1.93%  '
$ws.Range('E27').Value = '  -4.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.04'
$ws.Range('E28').Value = '  -2.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.40'
$ws.Range('E29').Value = '  -3.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.32'
$ws.Range('E30').Value = '  -2.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.56'
$ws.Range('E31').Value = '  -4.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.26'
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '582.10'
$ws.Range('E35').Value = '  -7.94%  '
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range('E37').Value = '  -0.78%  '
$ws.Range('E38').Value = '  +0.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.401'
$ws.Range('E39').Value = '  -2.71%  '
$ws.Range('E40').Value = '  +5.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.24'
$ws.Range('E41').Value = '  +19.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.49'
$ws.Range('E42').Value = '  -2.57%  '
$ws.Range('E43').Value = '  -6.72%  '
$ws.Range('D44').Value = '3.227.86'
$ws.Range('E44').Value = '  -2.54%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.07'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('E46').Value = '  -1.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.78'
$ws.Range('E47').Value = '  +7.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.40'
$ws.Range('E48').Value = '  +3.49%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  -1.86%  '
$ws.Range('E51').Value = '  -0.13%  '
